# Updated cryptos list on Wed Apr 12 08:28:20 UTC 2023 with GitHub Actions
#
# Applies the latest scrape refresh to the cryptocurrency table on the
# active worksheet: per-row Price (column D) and Volume(1h) (column E)
# updates, plus a few rows whose underlying coin shifted position in the
# source ranking (their Coin name / Link / Price / Volume all move
# together).
#
# Price values such as "1.002" are plain text in the sheet (inline
# strings), not numbers -- Excel would otherwise happily reinterpret a
# bare "1.002" as a numeric literal when assigned via .Value. Routing the
# write through a quick "@" (Text) number-format toggle forces the COM
# layer to keep it as a string, and resetting the cell style back to
# "Normal" afterwards avoids leaving any stray formatting behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Text
    $r.Style = "Normal"
}

Set-TextValue "D2" "30.034.95"
$ws.Range("E2").Value = "  -0.29%  "
Set-TextValue "D3" "1.873.08"
$ws.Range("E3").Value = "  -2.56%  "
Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "319.28"
$ws.Range("E5").Value = "  -3.59%  "
Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  +0.01%  "
Set-TextValue "D7" "0.5042"
$ws.Range("E7").Value = "  -3.44%  "
Set-TextValue "D8" "0.3949"
$ws.Range("E8").Value = "  -3.64%  "
Set-TextValue "D9" "0.08190"
$ws.Range("E9").Value = "  -4.02%  "
Set-TextValue "D10" "42.20"
$ws.Range("E10").Value = "  -2.91%  "
Set-TextValue "D11" "1.093"
$ws.Range("E11").Value = "  -3.14%  "
Set-TextValue "D12" "23.75"
$ws.Range("E12").Value = "  +5.81%  "
Set-TextValue "D13" "1.868.48"
$ws.Range("E13").Value = "  -2.70%  "
Set-TextValue "D14" "6.299"
$ws.Range("E14").Value = "  -2.06%  "
Set-TextValue "D15" "7.181"
$ws.Range("E15").Value = "  -3.40%  "
$ws.Range("E16").Value = "  +0.03%  "
Set-TextValue "D17" "92.00"
$ws.Range("E17").Value = "  -4.33%  "
Set-TextValue "D18" "0.00001090"
$ws.Range("E18").Value = "  -2.32%  "
Set-TextValue "D19" "0.06410"
$ws.Range("E19").Value = "  -4.43%  "
Set-TextValue "D20" "18.14"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("E21").Value = "  +0.03%  "
Set-TextValue "D22" "30.028.05"
$ws.Range("E22").Value = "  -0.37%  "
Set-TextValue "D23" "5.835"
$ws.Range("E23").Value = "  -3.47%  "
$ws.Range("E24").Value = "  -1.73%  "
Set-TextValue "D25" "2.173"
$ws.Range("E25").Value = "  -2.41%  "
Set-TextValue "D26" "2.086.64"
$ws.Range("E26").Value = "  -2.55%  "
Set-TextValue "D27" "21.47"
$ws.Range("E27").Value = "  +1.38%  "
Set-TextValue "D28" "160.32"
$ws.Range("E28").Value = "  +0.43%  "
Set-TextValue "D29" "2.226"
$ws.Range("E29").Value = "  -9.58%  "
Set-TextValue "D30" "127.34"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("E31").Value = "  -1.46%  "
Set-TextValue "D32" "0.1034"
$ws.Range("E32").Value = "  -2.32%  "
Set-TextValue "D33" "5.938"
$ws.Range("E33").Value = "  -3.05%  "
Set-TextValue "D34" "3.686"
$ws.Range("E34").Value = "  +1.22%  "
Set-TextValue "D35" "0.02429"
$ws.Range("E35").Value = "  -3.37%  "
Set-TextValue "D36" "5.236"
$ws.Range("E36").Value = "  +0.25%  "
Set-TextValue "D37" "0.06365"
$ws.Range("E37").Value = "  -3.62%  "
Set-TextValue "D38" "0.2146"
$ws.Range("E38").Value = "  -3.51%  "
Set-TextValue "D39" "1.176"
$ws.Range("E39").Value = "  -5.02%  "
Set-TextValue "D40" "8.503"
$ws.Range("E40").Value = "  -5.38%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D41" "1.219"
$ws.Range("E41").Value = "  -1.92%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D42" "0.6306"
$ws.Range("E42").Value = "  -3.57%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D43" "11.34"
$ws.Range("E43").Value = "  -3.63%  "
Set-TextValue "D44" "0.9991"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D45" "0.5918"
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "12.90"
$ws.Range("E46").Value = "  -3.27%  "
Set-TextValue "D47" "2.110"
$ws.Range("E47").Value = "  +1.13%  "
Set-TextValue "D48" "3.630"
$ws.Range("E48").Value = "  -3.89%  "
Set-TextValue "D49" "122.88"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D51" "77.48"
$ws.Range("E51").Value = "  -3.24%  "
